$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column BB (column 54) values for rows 1-82 ---
# Row 1 is the header/date row; rows 2-82 mirror the existing BA column (last
# forecast value repeated) except for rows 71-82 which carry updated
# (re-estimated) forecast figures from the refreshed model run.
$bbValues = @{
    1 = 45986;
    2 = 0.7596024993684409;
    3 = 5.398981879140436;
    4 = -5.469777829091811;
    5 = -0.5016308114541062;
    6 = -0.6214054032250829;
    7 = -1.41027180704431;
    8 = -0.7919227738455277;
    9 = 0.8281781605252121;
    10 = 0.6333580932638085;
    11 = -1.871721946158218;
    12 = 6.947074551031477;
    13 = -0.5538313209713124;
    14 = -2.926189405619382;
    15 = 9.498197854786355;
    16 = -0.4994768307962829;
    17 = -0.6902935135373696;
    18 = 1.660952585620919;
    19 = 0.5318249089792459;
    20 = 0.6320478207229741;
    21 = -0.1274490368921875;
    22 = -0.6653180828667331;
    23 = -2.768682804727675;
    24 = 2.972562358782;
    25 = 1.814509864363558;
    26 = 0.7290697448430592;
    27 = 4.434377037602189;
    28 = -3.830674548307485;
    29 = -0.3546604117736365;
    30 = 0.5841110884116603;
    31 = 0.3311457281599957;
    32 = -0.983338854947732;
    33 = 0.0367937104471423;
    34 = 2.049812641278834;
    35 = 2.087453263562125;
    36 = -1.438164479066856;
    37 = -0.02684958208732269;
    38 = 0.8237845450787233;
    39 = 2.86126545641612;
    40 = 0.5069532644027532;
    41 = -0.2174083661145261;
    42 = 0.1045745676356802;
    43 = 0.4936918743079417;
    44 = 1.322625689988016;
    45 = 0.6899772607160202;
    46 = 1.240850467000882;
    47 = 2.6;
    48 = -0.9;
    49 = 0.4;
    50 = 0.1;
    51 = 4.305581231918552;
    52 = -4.290231943195352;
    53 = -1.294763502326944;
    54 = 3.727166552773809;
    55 = -0.07597546477697392;
    56 = 1.620116340063873;
    57 = -3.620177172843626;
    58 = -0.9596300199833507;
    59 = 3.509382378526155;
    60 = -3.103464835402306;
    61 = -0.9738690852321383;
    62 = -3.206836580208446;
    63 = 1.5011358099934;
    64 = -0.5164324031961769;
    65 = -0.7267585812813877;
    66 = -1.863476264789497;
    67 = 0.5061198669764195;
    68 = -1.844723754442185;
    69 = -0.3113146714536583;
    70 = 0.8443262641745406;
    71 = 0.2255501838065186;
    72 = -1.137044400346582;
    73 = -0.4792091214565772;
    74 = 0.4580025317319599;
    75 = 0.1500769980497607;
    76 = 0.08692382407563778;
    77 = 0.1452293504052996;
    78 = 0.1409013763615072;
    79 = 0.1335834283383753;
    80 = 0.135872578757698;
    81 = 0.1363920206570818;
    82 = 0.1359471541599538
}

# Copy the column BA formatting (incl. header date style) into column BB
# first, then overwrite with the new values so styles line up with the rest
# of the sheet (BB1 gets the date style like BA1; BB2:BB82 stay unstyled like
# BA2:BA82).
$ws.Range("BA1:BA82").Copy()
$ws.Range("BB1:BB82").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 1; $r -le 82; $r++) {
    $ws.Cells.Item($r, 54).Value = $bbValues[$r]
}

# --- New row 83 (next period) ---
# Column A gets the next date stamp using the same date style as A82.
$ws.Range("A82").Copy()
$ws.Range("A83").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Cells.Item(83, 1).Value = 46934

# Column BB gets the new forecast value for row 83 (unstyled, like BA82).
$ws.Range("BA82").Copy()
$ws.Range("BB83").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Cells.Item(83, 54).Value = 0.1359748208526332
